$d = $word.ActiveDocument

$replacements = @(
    @("33×86=", "66×54="),
    @("93×20=", "81×56="),
    @("35×58=", "56×96="),
    @("66×20=", "51×82="),
    @("60×88=", "25×90="),
    @("80×45=", "13×93="),
    @("87×27=", "19×37="),
    @("29×68=", "32×23="),
    @("46×98=", "13×60="),
    @("49×57=", "18×42="),
    @("82×94=", "23×63="),
    @("44×18=", "95×70="),
    @("15×28=", "74×32="),
    @("11×49=", "89×20="),
    @("48×20=", "94×64="),
    @("28×82=", "90×93="),
    @("79×70=", "52×97="),
    @("32×45=", "88×12="),
    @("68×99=", "61×17="),
    @("82×13=", "30×73="),
    @("19×16=", "54×49="),
    @("56×75=", "11×53="),
    @("51×87=", "70×46="),
    @("36×28=", "29×75="),
    @("60×61=", "73×15=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
